# Uren registratie vrijdag en vandaag
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Friday of "Week 7" (row 23): hours logged for Zinedine, Robin, Carlo, Marc, Sam, Michiel
$ws.Range("D23:I23").Value = 4

# Monday of "Week 8" (row 27): hours logged for Rief, Zinedine, Robin, Carlo, Marc, Sam
$ws.Range("C27:H27").Value = 4

# Weekly "Totaal Game-Lab uren p/w" totals (manually entered, not formulas)
$ws.Range("B24").Value = 19
$ws.Range("B32").Value = 4

# Selection as last left by the author (Excel also scrolled the view to
# keep the new selection visible, topLeftCell A16)
$ws.Range("K28").Select()
